$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211588740348816
$ws.Range("B1").Value = 2.579480409622192
$ws.Range("C1").Value = 4.356475353240967
$ws.Range("D1").Value = 2.138078451156616
$ws.Range("E1").Value = 1.168679118156433
